$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - kept as text via leading apostrophe so the
# stored value matches the literal string from the source data feed.
$ws.Range("D2").Value = "'246.96"
$ws.Range("D3").Value = "'22.69"
$ws.Range("D4").Value = "'5.299"
$ws.Range("D5").Value = "'0.05723"
$ws.Range("D6").Value = "'3.439"
$ws.Range("D7").Value = "'0.8100"
$ws.Range("D8").Value = "'0.8682"
$ws.Range("D9").Value = "'0.1426"
$ws.Range("D10").Value = "'0.07356"
$ws.Range("D12").Value = "'0.03128"
$ws.Range("D13").Value = "'0.09401"
$ws.Range("D14").Value = "'3.917"
$ws.Range("D15").Value = "'0.001589"
$ws.Range("D16").Value = "'0.04824"
$ws.Range("D17").Value = "'0.0005843"
$ws.Range("D18").Value = "'0.006146"
$ws.Range("D19").Value = "'0.005123"
$ws.Range("D20").Value = "'0.0009971"
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("D22").Value = "'3.732"
$ws.Range("D23").Value = "'6.320"
$ws.Range("D24").Value = "'2.181"

$ws.Range("D42").Value = "'0.1067"
$ws.Range("D43").Value = "'0.002702"
$ws.Range("D44").Value = "'0.008159"
$ws.Range("D45").Value = "'0.00005621"
$ws.Range("D47").Value = "'0.6003"
$ws.Range("D48").Value = "'0.1767"
$ws.Range("D50").Value = "'0.01011"

# Volume(1h) (column E) text updates - "Bestin24h" marker moved from
# KickToken to CoinbaseStockToken.
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
